# Daily "remaining days" refresh.
#
# Column layout: A=row#, B=shop, C=address, D=total days ("总天"),
# E=remaining days ("剩余"), F=start date as an integer YYYYMMDD ("开始时间").
#
# For every data row, remaining = D - (today - F). The workbook had been
# generated with "today" = 2025-12-15..2025-12-22; this run advances the
# clock to 2025-12-23 and recomputes E for every row whose F parses as a
# clean 8-digit YYYYMMDD date. If the recomputed value would drop to 0 or
# below, the cycle is treated as renewed: E resets to the full D and F is
# bumped to the new "today" (2025-12-23). Rows with an unparseable F (not
# an 8-digit date) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "today" used to refresh the countdown.
$todayYear  = 2025
$todayMonth = 12
$todayDay   = 23
$todayInt   = $todayYear * 10000 + $todayMonth * 100 + $todayDay

function DaysFromCivil($y, $m, $d) {
    # Howard Hinnant's civil_from_days / days_from_civil algorithm.
    # Returns a monotonically increasing integer day-count usable for
    # plain subtraction between two dates - avoids relying on any
    # DateTime/TimeSpan support in the host.
    if ($m -le 2) {
        $y = $y - 1
    }
    if ($y -ge 0) {
        $era = [Math]::Floor($y / 400)
    } else {
        $era = [Math]::Floor(($y - 399) / 400)
    }
    $yoe = $y - $era * 400
    if ($m -gt 2) {
        $mp = $m - 3
    } else {
        $mp = $m + 9
    }
    $doy = [Math]::Floor((153 * $mp + 2) / 5) + $d - 1
    $doe = $yoe * 365 + [Math]::Floor($yoe / 4) - [Math]::Floor($yoe / 100) + $doy
    $result = $era * 146097 + $doe - 719468
    return $result
}

$todayDays = DaysFromCivil $todayYear $todayMonth $todayDay

$lastRow = $ws.Cells($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    $dVal = $dCell.Value2
    $fVal = $fCell.Value2

    if ($dVal -eq $null) { continue }
    if ($fVal -eq $null) { continue }

    $fStr = [string]([int64]$fVal)
    if ($fStr.Length -ne 8) {
        # Not a clean YYYYMMDD value (e.g. a data-entry typo) - skip,
        # exactly as the source data left it alone.
        continue
    }

    $fYear  = [int][Math]::Floor($fVal / 10000)
    $fMonth = [int][Math]::Floor(($fVal % 10000) / 100)
    $fDay   = [int]($fVal % 100)

    if ($fMonth -lt 1 -or $fMonth -gt 12) { continue }
    if ($fDay -lt 1 -or $fDay -gt 31) { continue }

    $fDays = DaysFromCivil $fYear $fMonth $fDay
    $elapsed = $todayDays - $fDays
    $newE = [int]$dVal - $elapsed

    if ($newE -le 0) {
        $eCell.Value = [int]$dVal
        $fCell.Value = $todayInt
    } else {
        $eCell.Value = $newE
    }
}

Write-Output "done"
